$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Solar" capacity values for 2023 (row 25) and 2024 (row 26)
$ws.Range("E25").Value = 409.874
$ws.Range("E26").Value = 330.385
